$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "68.863.20"
$ws.Range("E2").Value2 = "  -0.23%  "

$ws.Range("D3").Value2 = "3.849.62"
$ws.Range("E3").Value2 = "  +2.94%  "

$ws.Range("E4").Value2 = "  -0.03%  "

$ws.Range("D5").Value2 = "'601.17"
$ws.Range("E5").Value2 = "  -0.13%  "

$ws.Range("D6").Value2 = "'162.59"
$ws.Range("E6").Value2 = "  -3.03%  "

$ws.Range("D7").Value2 = "3.848.66"
$ws.Range("E7").Value2 = "  +2.95%  "

$ws.Range("E8").Value2 = "  +0.07%  "

$ws.Range("D9").Value2 = "'0.530"
$ws.Range("E9").Value2 = "  -1.44%  "

$ws.Range("E10").Value2 = "  -0.60%  "

$ws.Range("E11").Value2 = "  -2.41%  "

$ws.Range("E12").Value2 = "  -0.36%  "

$ws.Range("E13").Value2 = "  -3.27%  "

$ws.Range("E14").Value2 = "  -1.81%  "

$ws.Range("D15").Value2 = "4.497.91"
$ws.Range("E15").Value2 = "  +2.97%  "

$ws.Range("D16").Value2 = "3.860.54"
$ws.Range("E16").Value2 = "  +3.22%  "

$ws.Range("D17").Value2 = "69.025.17"
$ws.Range("E17").Value2 = "  +0.00%  "

$ws.Range("E18").Value2 = "  +2.58%  "

$ws.Range("E19").Value2 = "  -0.15%  "

$ws.Range("D20").Value2 = "'11.40"
$ws.Range("E20").Value2 = "  +4.19%  "

$ws.Range("D21").Value2 = "'17.10"
$ws.Range("E21").Value2 = "  -1.09%  "

$ws.Range("D22").Value2 = "'483.52"
$ws.Range("E22").Value2 = "  -1.79%  "

$ws.Range("D23").Value2 = "'0.719"
$ws.Range("E23").Value2 = "  -0.88%  "

$ws.Range("D25").Value2 = "'83.78"
$ws.Range("E25").Value2 = "  -1.18%  "

$ws.Range("E26").Value2 = "  -2.65%  "

$ws.Range("D27").Value2 = "'12.07"
$ws.Range("E27").Value2 = "  -2.24%  "

$ws.Range("D28").Value2 = "'10.01"
$ws.Range("E28").Value2 = "  -0.97%  "

$ws.Range("E29").Value2 = "  -0.10%  "

$ws.Range("E30").Value2 = "  -0.86%  "

$ws.Range("E31").Value2 = "  -1.33%  "

$ws.Range("D32").Value2 = "4.003.81"
$ws.Range("E32").Value2 = "  +3.03%  "

$ws.Range("E33").Value2 = "  -3.93%  "

$ws.Range("D34").Value2 = "'32.27"
$ws.Range("E34").Value2 = "  +2.45%  "

$ws.Range("D35").Value2 = "3.799.52"
$ws.Range("E35").Value2 = "  +3.32%  "

$ws.Range("E36").Value2 = "  -1.89%  "

$ws.Range("E37").Value2 = "  +1.44%  "

$ws.Range("E38").Value2 = "  +3.93%  "

$ws.Range("E39").Value2 = "  -0.80%  "

$ws.Range("D40").Value2 = "'1.00"

$ws.Range("D41").Value2 = "'0.318"
$ws.Range("E41").Value2 = "  -1.84%  "

$ws.Range("D42").Value2 = "'440.96"
$ws.Range("E42").Value2 = "  +2.04%  "

$ws.Range("E43").Value2 = "  -0.32%  "

$ws.Range("D44").Value2 = "'48.48"
$ws.Range("E44").Value2 = "  -0.34%  "

$ws.Range("E45").Value2 = "  -1.17%  "

$ws.Range("E47").Value2 = "  -1.15%  "

$ws.Range("D48").Value2 = "'26.55"
$ws.Range("E48").Value2 = "  +12.86%  "

$ws.Range("D49").Value2 = "'142.94"
$ws.Range("E49").Value2 = "  +1.28%  "

$ws.Range("D50").Value2 = "2.826.74"
$ws.Range("E50").Value2 = "  +1.55%  "

$ws.Range("D51").Value2 = "'0.0357"
$ws.Range("E51").Value2 = "  +1.45%  "
